# Update the "想去人数" (want-to-go count) column F values on both the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets to match the refreshed
# scrape output committed at 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 51
    6  = 43
    7  = 139
    9  = 20
    14 = 88
    15 = 463
    17 = 479
    18 = 408
    20 = 65
    23 = 1163
    24 = 2862
    27 = 548
    28 = 52
    29 = 1616
    31 = 454
    32 = 265
    33 = 391
    35 = 605
}

foreach ($sheetIndex in @(1, 4)) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
